# Scheduled-runner price refresh: update the recalculated market/profit
# columns (H-N) on a handful of rows across several sheets to match the
# latest pull. Values are plain numbers (no formulas are used in this
# workbook), so each touched cell is just overwritten with its new figure.
# A few cells that no longer have a computed value are cleared outright,
# and a couple of previously-empty "LeveProfitHQ" cells gain a new value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 856  # H33: was 977.4
$ws.Cells.Item(33, 9).Value = 729.1818  # I33: was 866.65216
$ws.Cells.Item(33, 11).Value = 729.1818  # K33: was 866.65216
$ws.Cells.Item(33, 13).Value = -500.1818  # M33: was -637.65216

$ws.Cells.Item(41, 8).Value = 1539.4286  # H41: was 1196.421
$ws.Cells.Item(41, 10).Value = 2574.25  # J41: was 1275.2222
$ws.Cells.Item(41, 12).Value = 2574.25  # L41: was 1275.2222
$ws.Cells.Item(41, 14).Value = -3454.25  # N41: was -2155.2222

$ws.Cells.Item(53, 8).Value = 439.6842  # H53: was 381.57144
$ws.Cells.Item(53, 10).Value = 906.3333  # J53: was 637.125
$ws.Cells.Item(53, 12).Value = 906.3333  # L53: was 637.125
$ws.Cells.Item(53, 14).Value = -2180.3333  # N53: was -1911.125

$ws.Cells.Item(99, 8).Value = 430.33334  # H99: was 482.46155
$ws.Cells.Item(99, 9).Value = 430.33334  # I99: was 497.83334
$ws.Cells.Item(99, 10).Value = 0  # J99: was 298
$ws.Cells.Item(99, 11).Value = 1291.00002  # K99: was 1493.50002
$ws.Cells.Item(99, 12).Value = 0  # L99: was 894
$ws.Cells.Item(99, 13).Value = 206.9999800000001  # M99: was 4.49998000000005
$ws.Cells.Item(99, 14).ClearContents()  # N99: was -3890

$ws.Cells.Item(137, 8).Value = 6870.1665  # H137: was 6457
$ws.Cells.Item(137, 9).Value = 4073.75  # I137: was 3558.8
$ws.Cells.Item(137, 11).Value = 12221.25  # K137: was 10676.4
$ws.Cells.Item(137, 13).Value = -9671.25  # M137: was -8126.400000000001

$ws.Cells.Item(138, 8).Value = 2522.1538  # H138: was 2522.738
$ws.Cells.Item(138, 10).Value = 3075.4827  # J138: was 3024.375
$ws.Cells.Item(138, 12).Value = 9226.4481  # L138: was 9073.125
$ws.Cells.Item(138, 14).Value = -19506.4481  # N138: was -19353.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11369358  # H32: was 11910725
$ws.Cells.Item(32, 9).Value = 12201179  # I32: was 12826847
$ws.Cells.Item(32, 11).Value = 12201179  # K32: was 12826847
$ws.Cells.Item(32, 13).Value = -12200892  # M32: was -12826560

$ws.Cells.Item(62, 8).Value = 20000  # H62: was 0
$ws.Cells.Item(62, 10).Value = 20000  # J62: was 0
$ws.Cells.Item(62, 12).Value = 20000  # L62: was 0
$ws.Cells.Item(62, 14).Value = -21248  # N62: was None

$ws.Cells.Item(65, 8).Value = 20000  # H65: was 0
$ws.Cells.Item(65, 10).Value = 20000  # J65: was 0
$ws.Cells.Item(65, 12).Value = 60000  # L65: was 0
$ws.Cells.Item(65, 14).Value = -66240  # N65: was None

$ws.Cells.Item(74, 8).Value = 10425356  # H74: was 11914503
$ws.Cells.Item(74, 9).Value = 20834852  # I74: was 22728838
$ws.Cells.Item(74, 10).Value = 15859  # J74: was 18733.1
$ws.Cells.Item(74, 11).Value = 20834852  # K74: was 22728838
$ws.Cells.Item(74, 12).Value = 15859  # L74: was 18733.1
$ws.Cells.Item(74, 13).Value = -20833978  # M74: was -22727964
$ws.Cells.Item(74, 14).Value = -17607  # N74: was -20481.1

$ws.Cells.Item(77, 8).Value = 10425356  # H77: was 11914503
$ws.Cells.Item(77, 9).Value = 20834852  # I77: was 22728838
$ws.Cells.Item(77, 10).Value = 15859  # J77: was 18733.1
$ws.Cells.Item(77, 11).Value = 104174260  # K77: was 113644190
$ws.Cells.Item(77, 12).Value = 79295  # L77: was 93665.5
$ws.Cells.Item(77, 13).Value = -104169892  # M77: was -113639822
$ws.Cells.Item(77, 14).Value = -88031  # N77: was -102401.5

$ws.Cells.Item(109, 8).Value = 59000  # H109: was 47000
$ws.Cells.Item(109, 10).Value = 59000  # J109: was 47000
$ws.Cells.Item(109, 12).Value = 59000  # L109: was 47000
$ws.Cells.Item(109, 14).Value = -61774  # N109: was -49774

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2800  # H86: was 2877.5557
$ws.Cells.Item(86, 9).Value = 2800  # I86: was 2499.6667
$ws.Cells.Item(86, 10).Value = 0  # J86: was 3633.3333
$ws.Cells.Item(86, 11).Value = 2800  # K86: was 2499.6667
$ws.Cells.Item(86, 12).Value = 0  # L86: was 3633.3333
$ws.Cells.Item(86, 13).Value = -1677  # M86: was -1376.6667
$ws.Cells.Item(86, 14).ClearContents()  # N86: was -5879.3333

$ws.Cells.Item(89, 8).Value = 2800  # H89: was 2877.5557
$ws.Cells.Item(89, 9).Value = 2800  # I89: was 2499.6667
$ws.Cells.Item(89, 10).Value = 0  # J89: was 3633.3333
$ws.Cells.Item(89, 11).Value = 14000  # K89: was 12498.3335
$ws.Cells.Item(89, 12).Value = 0  # L89: was 18166.6665
$ws.Cells.Item(89, 13).Value = -8384  # M89: was -6882.333500000001
$ws.Cells.Item(89, 14).ClearContents()  # N89: was -29398.6665

$ws.Cells.Item(94, 8).Value = 1007.8461  # H94: was 1013.1923
$ws.Cells.Item(94, 9).Value = 968.16  # I94: was 973.72
$ws.Cells.Item(94, 11).Value = 968.16  # K94: was 973.72
$ws.Cells.Item(94, 13).Value = -517.16  # M94: was -522.72

$ws.Cells.Item(107, 8).Value = 1972.75  # H107: was 1997.4445
$ws.Cells.Item(107, 9).Value = 1879.4  # I107: was 1897.7142
$ws.Cells.Item(107, 10).Value = 2128.3333  # J107: was 2346.5
$ws.Cells.Item(107, 11).Value = 1879.4  # K107: was 1897.7142
$ws.Cells.Item(107, 12).Value = 2128.3333  # L107: was 2346.5
$ws.Cells.Item(107, 13).Value = 40.59999999999991  # M107: was 22.28580000000011
$ws.Cells.Item(107, 14).Value = -5968.3333  # N107: was -6186.5

$ws.Cells.Item(134, 8).Value = 34474.367  # H134: was 36926.465
$ws.Cells.Item(134, 9).Value = 1216.6786  # I134: was 1256.5555
$ws.Cells.Item(134, 10).Value = 500082  # J134: was 1000014
$ws.Cells.Item(134, 11).Value = 3650.0358  # K134: was 3769.6665
$ws.Cells.Item(134, 12).Value = 1500246  # L134: was 3000042
$ws.Cells.Item(134, 13).Value = -1115.0358  # M134: was -1234.6665
$ws.Cells.Item(134, 14).Value = -1505316  # N134: was -3005112

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 557616.6  # H31: was 557677.9399999999
$ws.Cells.Item(31, 9).Value = 1777.3158  # I31: was 1802.4445
$ws.Cells.Item(31, 10).Value = 1016788.25  # J31: was 974584.5600000001
$ws.Cells.Item(31, 11).Value = 1777.3158  # K31: was 1802.4445
$ws.Cells.Item(31, 12).Value = 1016788.25  # L31: was 974584.5600000001
$ws.Cells.Item(31, 13).Value = -1482.3158  # M31: was -1507.4445
$ws.Cells.Item(31, 14).Value = -1017378.25  # N31: was -975174.5600000001

$ws.Cells.Item(34, 8).Value = 557616.6  # H34: was 557677.9399999999
$ws.Cells.Item(34, 9).Value = 1777.3158  # I34: was 1802.4445
$ws.Cells.Item(34, 10).Value = 1016788.25  # J34: was 974584.5600000001
$ws.Cells.Item(34, 11).Value = 1777.3158  # K34: was 1802.4445
$ws.Cells.Item(34, 12).Value = 1016788.25  # L34: was 974584.5600000001
$ws.Cells.Item(34, 13).Value = -1575.3158  # M34: was -1600.4445
$ws.Cells.Item(34, 14).Value = -1017192.25  # N34: was -974988.5600000001

$ws.Cells.Item(58, 8).Value = 2347.261  # H58: was 2214.3215
$ws.Cells.Item(58, 9).Value = 2109.611  # I58: was 1999.4348
$ws.Cells.Item(58, 11).Value = 2109.611  # K58: was 1999.4348
$ws.Cells.Item(58, 13).Value = -1906.611  # M58: was -1796.4348

$ws.Cells.Item(107, 8).Value = 3004.3333  # H107: was 2843.8333
$ws.Cells.Item(107, 9).Value = 3000  # I107: was 2810
$ws.Cells.Item(107, 11).Value = 3000  # K107: was 2810
$ws.Cells.Item(107, 13).Value = -1080  # M107: was -890

$ws.Cells.Item(136, 8).Value = 2347.261  # H136: was 2214.3215
$ws.Cells.Item(136, 9).Value = 2109.611  # I136: was 1999.4348
$ws.Cells.Item(136, 11).Value = 6328.833  # K136: was 5998.3044
$ws.Cells.Item(136, 13).Value = -3778.833  # M136: was -3448.3044

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(105, 8).Value = 9625  # H105: was 9666.666999999999
$ws.Cells.Item(105, 10).Value = 9625  # J105: was 9666.666999999999
$ws.Cells.Item(105, 12).Value = 28875  # L105: was 29000.001
$ws.Cells.Item(105, 14).Value = -34117  # N105: was -34242.001

$ws.Cells.Item(113, 8).Value = 1129.2858  # H113: was 1002.82355
$ws.Cells.Item(113, 9).Value = 466.8  # I113: was 524.8
$ws.Cells.Item(113, 10).Value = 1497.3334  # J113: was 1202
$ws.Cells.Item(113, 11).Value = 1400.4  # K113: was 1574.4
$ws.Cells.Item(113, 12).Value = 4492.0002  # L113: was 3606
$ws.Cells.Item(113, 13).Value = 769.5999999999999  # M113: was 595.6000000000001
$ws.Cells.Item(113, 14).Value = -8832.0002  # N113: was -7946

$ws.Cells.Item(122, 8).Value = 1436.2  # H122: was 1472.5385
$ws.Cells.Item(122, 9).Value = 715.3333  # I122: was 778.4
$ws.Cells.Item(122, 10).Value = 1916.7778  # J122: was 1906.375
$ws.Cells.Item(122, 11).Value = 6437.9997  # K122: was 7005.599999999999
$ws.Cells.Item(122, 12).Value = 17251.0002  # L122: was 17157.375
$ws.Cells.Item(122, 13).Value = -3987.9997  # M122: was -4555.599999999999
$ws.Cells.Item(122, 14).Value = -22151.0002  # N122: was -22057.375

$ws.Cells.Item(123, 8).Value = 54979.633  # H123: was 57933.61
$ws.Cells.Item(123, 9).Value = 3403.75  # I123: was 3935.6667
$ws.Cells.Item(123, 11).Value = 10211.25  # K123: was 11807.0001
$ws.Cells.Item(123, 13).Value = -7761.25  # M123: was -9357.000100000001

$ws.Cells.Item(124, 8).Value = 1575.7142  # H124: was 1675.7142
$ws.Cells.Item(124, 9).Value = 1575.7142  # I124: was 1675.7142
$ws.Cells.Item(124, 11).Value = 4727.142599999999  # K124: was 5027.142599999999
$ws.Cells.Item(124, 13).Value = 182.8574000000008  # M124: was -117.1425999999992

$ws.Cells.Item(125, 8).Value = 3500  # H125: was 10892.556
$ws.Cells.Item(125, 10).Value = 4000  # J125: was 13433.286
$ws.Cells.Item(125, 12).Value = 12000  # L125: was 40299.858
$ws.Cells.Item(125, 14).Value = -21840  # N125: was -50139.858

$ws.Cells.Item(129, 8).Value = 83335590  # H129: was 1048.2858
$ws.Cells.Item(129, 9).Value = 2000  # I129: was 556.3333
$ws.Cells.Item(129, 10).Value = 111113460  # J129: was 4000
$ws.Cells.Item(129, 11).Value = 6000  # K129: was 1668.9999
$ws.Cells.Item(129, 12).Value = 333340380  # L129: was 12000
$ws.Cells.Item(129, 13).Value = -1000  # M129: was 3331.0001
$ws.Cells.Item(129, 14).Value = -333350380  # N129: was -22000

$ws.Cells.Item(130, 8).Value = 5000  # H130: was 4810.9375
$ws.Cells.Item(130, 9).Value = 0  # I130: was 1975
$ws.Cells.Item(130, 11).Value = 0  # K130: was 5925
$ws.Cells.Item(130, 13).ClearContents()  # M130: was -905

$ws.Cells.Item(132, 8).Value = 2407.2778  # H132: was 2050.0454
$ws.Cells.Item(132, 9).Value = 2493.4  # I132: was 2144.5
$ws.Cells.Item(132, 10).Value = 2299.625  # J132: was 1936.7
$ws.Cells.Item(132, 11).Value = 22440.6  # K132: was 19300.5
$ws.Cells.Item(132, 12).Value = 20696.625  # L132: was 17430.3
$ws.Cells.Item(132, 13).Value = -19910.6  # M132: was -16770.5
$ws.Cells.Item(132, 14).Value = -25756.625  # N132: was -22490.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3674  # H113: was 3782.2222
$ws.Cells.Item(113, 9).Value = 3623.75  # I113: was 3755.7144
$ws.Cells.Item(113, 11).Value = 3623.75  # K113: was 3755.7144
$ws.Cells.Item(113, 13).Value = -1453.75  # M113: was -1585.7144

$ws.Cells.Item(132, 8).Value = 100003280  # H132: was 111114510
$ws.Cells.Item(132, 9).Value = 125003000  # I132: was 142860260
$ws.Cells.Item(132, 11).Value = 375009000  # K132: was 428580780
$ws.Cells.Item(132, 13).Value = -375006470  # M132: was -428578250

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 54874.8  # H7: was 57610.105
$ws.Cells.Item(7, 9).Value = 4075.1538  # I7: was 4172.75
$ws.Cells.Item(7, 11).Value = 4075.1538  # K7: was 4172.75
$ws.Cells.Item(7, 13).Value = -3963.1538  # M7: was -4060.75

$ws.Cells.Item(46, 8).Value = 2553.3948  # H46: was 2496.9
$ws.Cells.Item(46, 9).Value = 2465.3572  # I46: was 2395.9
$ws.Cells.Item(46, 11).Value = 2465.3572  # K46: was 2395.9
$ws.Cells.Item(46, 13).Value = -2277.3572  # M46: was -2207.9

$ws.Cells.Item(122, 8).Value = 5478.6816  # H122: was 5412.2173
$ws.Cells.Item(122, 9).Value = 4428.8  # I122: was 4398.875
$ws.Cells.Item(122, 11).Value = 13286.4  # K122: was 13196.625
$ws.Cells.Item(122, 13).Value = -10836.4  # M122: was -10746.625

$ws.Cells.Item(126, 8).Value = 54874.8  # H126: was 57610.105
$ws.Cells.Item(126, 9).Value = 4075.1538  # I126: was 4172.75
$ws.Cells.Item(126, 11).Value = 12225.4614  # K126: was 12518.25
$ws.Cells.Item(126, 13).Value = -9755.4614  # M126: was -10048.25

$ws.Cells.Item(132, 8).Value = 184819.64  # H132: was 225068.56
$ws.Cells.Item(132, 9).Value = 169002.67  # I132: was 202222.4
$ws.Cells.Item(132, 10).Value = 203800  # J132: was 253626.25
$ws.Cells.Item(132, 11).Value = 507008.01  # K132: was 606667.2
$ws.Cells.Item(132, 12).Value = 611400  # L132: was 760878.75
$ws.Cells.Item(132, 13).Value = -504478.01  # M132: was -604137.2
$ws.Cells.Item(132, 14).Value = -616460  # N132: was -765938.75

$ws.Cells.Item(136, 8).Value = 41141.066  # H136: was 48284.152
$ws.Cells.Item(136, 9).Value = 3896.3333  # I136: was 3927.5334
$ws.Cells.Item(136, 10).Value = 92710.69500000001  # J136: was 108770.45
$ws.Cells.Item(136, 11).Value = 11688.9999  # K136: was 11782.6002
$ws.Cells.Item(136, 12).Value = 278132.085  # L136: was 326311.35
$ws.Cells.Item(136, 13).Value = -9138.999899999999  # M136: was -9232.600199999999
$ws.Cells.Item(136, 14).Value = -283232.085  # N136: was -331411.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 40006268  # H62: was 50005892
$ws.Cells.Item(62, 10).Value = 50005864  # J62: was 66671892
$ws.Cells.Item(62, 12).Value = 50005864  # L62: was 66671892
$ws.Cells.Item(62, 14).Value = -50007112  # N62: was -66673140

$ws.Cells.Item(65, 8).Value = 40006268  # H65: was 50005892
$ws.Cells.Item(65, 10).Value = 50005864  # J65: was 66671892
$ws.Cells.Item(65, 12).Value = 250029320  # L65: was 333359460
$ws.Cells.Item(65, 14).Value = -250035560  # N65: was -333365700

$ws.Cells.Item(107, 8).Value = 18519952  # H107: was 19232202
$ws.Cells.Item(107, 9).Value = 26317380  # I107: was 27779376
$ws.Cells.Item(107, 11).Value = 78952140  # K107: was 83338128
$ws.Cells.Item(107, 13).Value = -78950220  # M107: was -83336208

$ws.Cells.Item(113, 8).Value = 1143.7  # H113: was 1333.375
$ws.Cells.Item(113, 9).Value = 1104.1111  # I113: was 1309.5714
$ws.Cells.Item(113, 11).Value = 3312.3333  # K113: was 3928.7142
$ws.Cells.Item(113, 13).Value = -1142.3333  # M113: was -1758.7142
